# Apply cryptos list update (commit: "Updated cryptos list on Wed Aug 23 08:52:34 UTC 2023
# with GitHub Actions") -- refreshed Price/Volume(1h) figures for the existing coins, and
# the table dropped BabyDogeCoin while gaining Aptos at the bottom, shifting rows 46-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds values such as "26.078.61" or "0.5222" that are stored as
# plain text in the source data (note the multi-dot thousands format on BTC/ETH-scale
# rows), not numbers. Force text format first so Excel does not reinterpret strings like
# "215.56" as a number, then restore the default (unstyled) look once the text is in place.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value2 = '26.078.61'
$ws.Range('E2').Value2 = '  -0.31%  '
$ws.Range('D3').Value2 = '1.646.10'
$ws.Range('E4').Value2 = '  -0.17%  '
$ws.Range('D5').Value2 = '215.56'
$ws.Range('E5').Value2 = '  +2.24%  '
$ws.Range('D6').Value2 = '0.5222'
$ws.Range('E6').Value2 = '  +0.41%  '
$ws.Range('E7').Value2 = '  -0.18%  '
$ws.Range('D9').Value2 = '0.06379'
$ws.Range('E9').Value2 = '  +0.85%  '
$ws.Range('D10').Value2 = '20.85'
$ws.Range('E10').Value2 = '  -1.61%  '
$ws.Range('D11').Value2 = '0.07664'
$ws.Range('E11').Value2 = '  +1.52%  '
$ws.Range('D12').Value2 = '1.652.44'
$ws.Range('E12').Value2 = '  -1.07%  '
$ws.Range('D13').Value2 = '4.424'
$ws.Range('E13').Value2 = '  -0.43%  '
$ws.Range('D14').Value2 = '1.868.71'
$ws.Range('E14').Value2 = '  -1.61%  '
$ws.Range('D15').Value2 = '0.5550'
$ws.Range('E15').Value2 = '  +1.02%  '
$ws.Range('D16').Value2 = '0.0₅8261'
$ws.Range('E16').Value2 = '  +2.99%  '
$ws.Range('D17').Value2 = '65.01'
$ws.Range('D18').Value2 = '26.089.97'
$ws.Range('E18').Value2 = '  -0.28%  '
$ws.Range('E19').Value2 = '  -0.11%  '
$ws.Range('D20').Value2 = '4.734'
$ws.Range('E20').Value2 = '  -0.58%  '
$ws.Range('D21').Value2 = '188.27'
$ws.Range('E21').Value2 = '  +0.66%  '
$ws.Range('E22').Value2 = '  -1.11%  '
$ws.Range('D23').Value2 = '6.232'
$ws.Range('E23').Value2 = '  +0.31%  '
$ws.Range('E24').Value2 = '  -0.17%  '
$ws.Range('D25').Value2 = '146.16'
$ws.Range('E25').Value2 = '  -2.66%  '
$ws.Range('D26').Value2 = '0.1220'
$ws.Range('E26').Value2 = '  -1.62%  '
$ws.Range('D27').Value2 = '7.427'
$ws.Range('E27').Value2 = '  -0.89%  '
$ws.Range('D28').Value2 = '15.84'
$ws.Range('E28').Value2 = '  +0.12%  '
$ws.Range('D29').Value2 = '1.381'
$ws.Range('E29').Value2 = '  +2.36%  '
$ws.Range('D30').Value2 = '0.05953'
$ws.Range('E30').Value2 = '  -5.49%  '
$ws.Range('D31').Value2 = '1.266'
$ws.Range('E31').Value2 = '  -1.39%  '
$ws.Range('D33').Value2 = '3.405'
$ws.Range('E33').Value2 = '  -3.38%  '
$ws.Range('D34').Value2 = '1.665'
$ws.Range('E34').Value2 = '  +1.15%  '
$ws.Range('D35').Value2 = '0.9947'
$ws.Range('E35').Value2 = '  -1.12%  '
$ws.Range('D36').Value2 = '2.394'
$ws.Range('E36').Value2 = '  -0.43%  '
$ws.Range('D37').Value2 = '2.752'
$ws.Range('E37').Value2 = '  -0.27%  '
$ws.Range('D38').Value2 = '0.5643'
$ws.Range('E38').Value2 = '  -6.86%  '
$ws.Range('E39').Value2 = '  +0.23%  '
$ws.Range('D40').Value2 = '0.8600'
$ws.Range('E40').Value2 = '  -0.61%  '
$ws.Range('D41').Value2 = '5.841'
$ws.Range('E41').Value2 = '  -4.47%  '
$ws.Range('E42').Value2 = '  -0.28%  '
$ws.Range('D43').Value2 = '1.029.36'
$ws.Range('E43').Value2 = '  -7.51%  '
$ws.Range('D44').Value2 = '99.97'
$ws.Range('E44').Value2 = '  -0.60%  '
$ws.Range('D45').Value2 = '1.795.42'
$ws.Range('E45').Value2 = '  -1.47%  '
$ws.Range('B46').Value2 = 'Aave'
$ws.Range('C46').Value2 = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value2 = '55.89'
$ws.Range('E46').Value2 = '  +0.55%  '
$ws.Range('B47').Value2 = 'EnergySwap'
$ws.Range('C47').Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value2 = '8.101'
$ws.Range('E47').Value2 = '  +0.13%  '
$ws.Range('D48').Value2 = '0.9989'
$ws.Range('E48').Value2 = '  +0.00%  '
$ws.Range('B49').Value2 = 'Cronos'
$ws.Range('C49').Value2 = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value2 = '0.05164'
$ws.Range('E49').Value2 = '  -1.41%  '
$ws.Range('B50').Value2 = 'Mantle'
$ws.Range('C50').Value2 = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value2 = '0.4220'
$ws.Range('E50').Value2 = '  -0.49%  '
$ws.Range('B51').Value2 = 'Aptos'
$ws.Range('C51').Value2 = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').Value2 = '5.946'
$ws.Range('E51').Value2 = '  +0.28%  '

$priceRange.Style = "Normal"
